# Apply new TPM values and remove obsolete rows for Cxcl2-Cxcr1 LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed values for the remaining sending-cluster rows (ECs, FAPs, MuSCs, Resolving-Mac)
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 1.400677
$ws.Range("H2").Value = 4.202031
$ws.Range("I2").Value = 0.009106128952548741
$ws.Range("J2").Value = 0.009106128952548741
$ws.Range("M2").Value = 0.009727666666666667
$ws.Range("N2").Value = 0.029183
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01362531896366667
$ws.Range("R2").Value = 0.122627870673
$ws.Range("S2").Value = 0.009106128952548741
$ws.Range("T2").Value = 0.009106128952548741
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 6.974902666666668
$ws.Range("H3").Value = 20.924708
$ws.Range("I3").Value = 0.04534547444852937
$ws.Range("J3").Value = 0.04534547444852936
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.009727666666666667
$ws.Range("N3").Value = 0.029183
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.06784952817377779
$ws.Range("R3").Value = 0.6106457535640001
$ws.Range("S3").Value = 0.04534547444852937
$ws.Range("T3").Value = 0.04534547444852936
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 1.381819666666667
$ws.Range("H4").Value = 4.145459
$ws.Range("I4").Value = 0.00898353301570211
$ws.Range("J4").Value = 0.00898353301570211
$ws.Range("M4").Value = 0.009727666666666667
$ws.Range("N4").Value = 0.029183
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01344188111077778
$ws.Range("R4").Value = 0.120976929997
$ws.Range("S4").Value = 0.00898353301570211
$ws.Range("T4").Value = 0.00898353301570211
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 144.0595526666667
$ws.Range("H5").Value = 432.178658
$ws.Range("I5").Value = 0.9365648635832199
$ws.Range("J5").Value = 0.9365648635832198
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009727666666666667
$ws.Range("N5").Value = 0.029183
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1.401363308490444
$ws.Range("R5").Value = 12.612269776414
$ws.Range("S5").Value = 0.9365648635832199
$ws.Range("T5").Value = 0.9365648635832198

# Remove the rows for cluster pairs that no longer apply (old rows 6-9)
$ws.Rows("6:9").Delete()
